# Adds 2022-Q3 data:
#  - "总计" summary sheet gets a new row for 2022-Q3 (inserted above the
#    existing 2022-Q2 row)
#  - a new "2022-Q3" worksheet is inserted right before the existing
#    "2022-Q2" worksheet, holding the per-fund holdings for the quarter

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. "总计" sheet: push the existing 2022-Q2 row down and insert the new
#    2022-Q3 row above it.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("D3").Value = 0.27

# row-insert leaves A2 without its index-column style and B2:D2 with a
# stray border style picked up from the row above - line both rows back
# up with the existing formatting (A3 / B3:D3 already carry the correct,
# original styling).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B3:D3").Copy()
$wsTotal.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Insert a brand new "2022-Q3" sheet right before the existing
#    "2022-Q2" sheet and fill it with the quarter's fund holdings.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

# Header formatting (bold + thin border) matches the "总计" sheet's style.
$wsTotal.Range("B1:D1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Index column (A) formatting.
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns B, D:G hold text (not numbers) even though several look numeric
# (e.g. the fund code has a leading zero, percentages are stored as plain
# text) - force text formatting so Excel doesn't coerce them to numbers.
# Column C (fund name) is already non-numeric text, no coercion risk.
$wsQ3.Range("B2:B3").NumberFormat = "@"
$wsQ3.Range("D2:G3").NumberFormat = "@"

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "001735"
$wsQ3.Range("C2").Value = "广发百发大数据策略成长灵活配置混合E"
$wsQ3.Range("D2").Value = "0.51"
$wsQ3.Range("E2").Value = "43.77"
$wsQ3.Range("F2").Value = "1.63"
$wsQ3.Range("G2").Value = "0.0083"
$wsQ3.Range("H2").Value = 3

$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "001734"
$wsQ3.Range("C3").Value = "广发百发大数据策略成长灵活配置混合A"
$wsQ3.Range("D3").Value = "0.27"
$wsQ3.Range("E3").Value = "43.77"
$wsQ3.Range("F3").Value = "1.63"
$wsQ3.Range("G3").Value = "0.0044"
$wsQ3.Range("H3").Value = 3

# Keep "2022-Q2" as the active tab, matching the original workbook.
$wb.Worksheets.Item(3).Activate()
